# 4.7.1 lab workbook update
# Adds the "+37-G{row}/2" correction term (matching the N/O-column pattern
# used earlier in the sheet) to the E/F angle-difference formulas for rows
# 35-62, and updates the sheet's view/zoom/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: standalone formulas (outside the shared-formula range) -----
$ws.Range("E35").Formula = "=H35-180+37-G35/2"
$ws.Range("F35").Formula = "=I35-180+37-G35/2"

# --- Column E, rows 36-62: shared formula anchored at E36 ---------------
$ws.Range("E36:E62").Formula = "=H36-180+37-G36/2"

# --- Column F, rows 36-62: update each row's formula ---------------------
For ($r = 36; $r -le 62; $r++) {
    $ws.Range("F$r").Formula = "=I$r-180+37-G$r/2"
}

# --- View state: zoom + scroll position + selection ----------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("C51").Select()
